# Remove the old Jekyll-site footer boilerplate that used to trail every
# "Requisitos" block: the blank spacer paragraph, the "Ver no Jupiter..."
# links line, and the "(c) 2020 ... Creative Commons Attribution" notice.
# The paragraph with the last course requisite (and everything that
# follows the footer, e.g. the trailing page-break paragraph) is left
# untouched.

$d = $word.ActiveDocument

# Locate the "Ver no Jupiter..." paragraph by its text.
$verFind = $d.Content
$verFound = $verFind.Find.Execute("Ver no Jupiter Salvar em pdf Salvar em docx", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

# Locate the copyright paragraph via a fragment that avoids the special
# "(c)" glyph (keeps the Find call encoding-agnostic).
$copyrightFind = $d.Content
$copyrightFound = $copyrightFind.Find.Execute("Powered by Jekyll and Github pages", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($verFound -and $copyrightFound) {
    # Translate the found character offsets into 1-based paragraph indexes
    # (Range.Paragraphs.Count over [0, pos+1) counts every paragraph that
    # the point at "pos" touches, including a partial one).
    $verIdx = $d.Range(0, $verFind.Start + 1).Paragraphs.Count
    $copyrightIdx = $d.Range(0, $copyrightFind.Start + 1).Paragraphs.Count

    # The blank spacer paragraph sits immediately before "Ver no Jupiter...".
    $startPara = $d.Paragraphs.Item($verIdx - 1)
    $endPara = $d.Paragraphs.Item($copyrightIdx)

    $deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
    $deleteRange.Delete()
}
